$d = $word.ActiveDocument

# --- Remove the stray _GoBack bookmark from the "created a separate reserved
#     word file..." paragraph (it will be re-added at the end of the new
#     content below). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Append the new journal entry:
#       (empty paragraph)
#       #date:10-6-19
#       <tab>finished with the removal of left recursion still need to
#            create the first and follows for the grammar.   [+ _GoBack bookmark]
#     Word stores the tab character as its own <w:tab/> run, so we build the
#     new paragraphs from raw WordprocessingML via Range.InsertXML rather
#     than Range.Text (which would serialize the tab as a literal character). ---
$lastPara = $d.Paragraphs.Last
$anchor = $lastPara.Range
$anchor.InsertParagraphAfter()

$newRange = $d.Paragraphs.Last.Range

$openXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:r>
              <w:t>#date:10-6-19</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:tab/>
              <w:t>finished with the removal of left recursion still need to create the first and follows for the grammar.</w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
          <w:sectPr/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newRange.InsertXML($openXml)
